$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "VelocidadePolPorMinuto" column (C) for all machine rows to a
# standard value of 1000.
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 1000

# Replace the old "CAL-24x5" calendar reference with the new standard
# calendar "CAL-PADRAO-5x8" for all machine rows (column H - CalendarioId).
$ws.Range("H2").Value = "CAL-PADRAO-5x8"
$ws.Range("H3").Value = "CAL-PADRAO-5x8"
$ws.Range("H4").Value = "CAL-PADRAO-5x8"

# Move selection to H4, matching the last edited cell before changing the grid.
$ws.Range("H4").Select()
